# Update gh-pages output (想去人数 column F) for the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 9069
$ws1.Range("F3").Value = 203
$ws1.Range("F4").Value = 464
$ws1.Range("F5").Value = 453

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 9069
$ws4.Range("F3").Value = 203
$ws4.Range("F4").Value = 464
$ws4.Range("F6").Value = 453
